$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("I5").Value = "Compound sentence; S (≤30 characters); Accuracy validation"

# Row 35
$ws.Range("E35").Value = "mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama mama "
$ws.Range("F35").Value = "මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම මම "

# Row 37
$ws.Range("C37").Value = "M"
$ws.Range("D37").Value = " mama pansal yannee naee "
$ws.Range("E37").Value = "mama pansal yannee naee"
$ws.Range("F37").Value = "මම පන්සල් යන්නේ නෑ"
$ws.Range("I37").Value = "Daily language usage; Sentence structure; M (≥30 characters); Real-time output update behavior"

# Row 38
$ws.Range("C38").Value = "M"
$ws.Range("D38").Value = "mama gedhara yannee naee , oyaata kohomadha? "
$ws.Range("E38").Value = "mama gedhara yannee naee , oyaata kohomadha?"
$ws.Range("F38").Value = "මම ගෙදර යන්නේ නෑ , ඔයාට කොහොමද?"
$ws.Range("I38").Value = "Usability flow; M (≥30 characters); Real-time output update behavior"
